# Update Name of Algo
# Apply targeted numeric corrections to result_data_KNN sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = -12.466
$ws.Range("B12").Value = 4.935
$ws.Range("C23").Value = -12.748
$ws.Range("B27").Value = 5.415000000000001
$ws.Range("C28").Value = -12.809
$ws.Range("B32").Value = 6.455
$ws.Range("C32").Value = -12.93
$ws.Range("C34").Value = -11.573
$ws.Range("B36").Value = 8.705
$ws.Range("B38").Value = 6.031
$ws.Range("C42").Value = -12.221
$ws.Range("B46").Value = 6.311
$ws.Range("C49").Value = -12.827
$ws.Range("B54").Value = 5.741
$ws.Range("C54").Value = -13.017
$ws.Range("B55").Value = 4.726
$ws.Range("B56").Value = 4.403
$ws.Range("B67").Value = 5.516
$ws.Range("B69").Value = 5.367
$ws.Range("B72").Value = 5.697
$ws.Range("C78").Value = -12.291
$ws.Range("C80").Value = -11.206
$ws.Range("B83").Value = 5.883
$ws.Range("B86").Value = 5.276
$ws.Range("B91").Value = 5.685
$ws.Range("B93").Value = 5.423999999999999
$ws.Range("C97").Value = -11.242
$ws.Range("B99").Value = 5.109999999999999
$ws.Range("C99").Value = -12.225
$ws.Range("C101").Value = -12.608
$ws.Range("B104").Value = 8.57
